$wb = $excel.ActiveWorkbook

# The same set of rows (F2, F3, F8, F9, F10 = "想去人数") are duplicated on both
# the "展览" sheet and the "全部类型" sheet; both need to be updated identically.
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 616
    $ws.Range("F3").Value = 469
    $ws.Range("F8").Value = 1080
    $ws.Range("F9").Value = 3871
    $ws.Range("F10").Value = 78
}
